$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52, shifting existing rows 52-78 down to 53-79
$ws.Rows.Item(52).Insert()

# Populate the new row 52 with the new data record
$ws.Cells.Item(52, 1).Value = 7
$ws.Cells.Item(52, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(52, 3).Value = "Ñuble"
$ws.Cells.Item(52, 4).Value = 45068
$ws.Cells.Item(52, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(52, 5).Value = 16
$ws.Cells.Item(52, 6).Value = 100112001
$ws.Cells.Item(52, 7).Value = "Berenjena"
$ws.Cells.Item(52, 8).Value = "Sin especificar"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 90
$ws.Cells.Item(52, 11).Value = 6000
$ws.Cells.Item(52, 12).Value = 7000
$ws.Cells.Item(52, 13).Value = 6556
$ws.Cells.Item(52, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(52, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(52, 16).Value = 109
$ws.Cells.Item(52, 17).Value = 60
$ws.Cells.Item(52, 18).Value = "Hortaliza"
